$wb = $excel.ActiveWorkbook

# --- Sheet "OFF": update row 3 (R) values ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 308
$wsOff.Range("C3").Value = 220
$wsOff.Range("D3").Value = 72
$wsOff.Range("E3").Value = 36

# --- Sheet "DEF": update row 3 (R) values ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 312
$wsDef.Range("C3").Value = 215
$wsDef.Range("D3").Value = 61
$wsDef.Range("E3").Value = 29
$wsDef.Range("F3").Value = 5
